$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (2008年 and 2009年), shifting rows 4/5 up to become 2/3
$ws.Rows("2:3").Delete()
